$wb = $excel.ActiveWorkbook

function Update-WantToGo($Sheet, $RowValues) {
    foreach ($row in $RowValues.Keys) {
        $Sheet.Range("F$row").Value = $RowValues[$row]
    }
}

# Sheet "展览" (Exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
Update-WantToGo $wsExhibit @{
    2  = 310
    3  = 13880
    6  = 182
    7  = 279
    9  = 12
    13 = 52
    14 = 455
    15 = 5851
    16 = 138
    17 = 92
    18 = 980
    19 = 114
    20 = 59
    21 = 155
    22 = 268
}

# Sheet "全部类型" (All types)
$wsAll = $wb.Worksheets.Item("全部类型")
Update-WantToGo $wsAll @{
    2  = 310
    3  = 13880
    6  = 182
    7  = 279
    9  = 12
    13 = 52
    14 = 455
    15 = 5851
    16 = 138
    17 = 92
    18 = 980
    19 = 114
    20 = 59
    21 = 156
    22 = 268
}
